$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("neg_reaction19")

$values = 1,2,3,5,6,7,8,9,10,13,14,15,16,18,19,20,22,23,24,25,28,29,30,31,32,33,34,35,36,38,42,43,45,46,47,48,50,51,52,53,55,56,58,59,60,61,62,63,64,65,66,67,68,69,71,72,73,74

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $values[$i]
}
